$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.47%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.48%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.761"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.26%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06071"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.74%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.718"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.10%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8662"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.28%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.47%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05008"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.62%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07121"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.36%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03062"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.42%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09118"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.22%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001536"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.54%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006088"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-94.18%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006194"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.57%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.449"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.51%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.168"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.61%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.25%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.37%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.62%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.091"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.26%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04259"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.14%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.42%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003912"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-8.93%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.05%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-18.86%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03884"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.40%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004133"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-33.97%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01495"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.11%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002209"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.41%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005348"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.17%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.05%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "7.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1353"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-46.51%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.05%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
